# "added week 1 day 1 booth"
# Add a Brand lookup column (C) to the VLookup worksheet, driven by a
# second VLOOKUP against the ID/Brand/Product table in E:G, and move the
# active selection the way the author left it (VLookup!H15 / HLookup!E17,
# with VLookup as the active tab).

$wb = $excel.ActiveWorkbook
$wsV = $wb.Worksheets.Item("VLookup")
$wsH = $wb.Worksheets.Item("HLookup")

# Header for the new column.
$wsV.Range("C1").Value = "Brand"

# Row 2 is a standalone formula (matches B2's pattern on this sheet).
$wsV.Range("C2").Formula = "=VLOOKUP(A2,`$E`$3:`$G`$7,2, FALSE)"

# Rows 3-11 share one formula (matches B3:B11's shared-formula pattern).
$wsV.Range("C3:C11").Formula = "=VLOOKUP(A3,`$E`$3:`$G`$7,2, FALSE)"

# Match column A's (unstyled/default) look for the new shared-formula cells.
$wsV.Range("A3:A11").Copy()
$wsV.Range("C3:C11").PasteSpecial(-4122)
$wsV.Range("C3:C11").Borders.LineStyle = -4142

# Restore the selections/active sheet the workbook was left on.
$wsH.Range("E17").Select() | Out-Null
$wsV.Activate() | Out-Null
$wsV.Range("H15").Select() | Out-Null
